# [LC-850] Release of LetsCo OS v1.3.0
# Rename generic "GPn"/"BPn" KPI name codes to two-digit "GPnn"/"BPnn" codes,
# and reset the sheet's scroll/selection state back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename KPI name codes (column B) ---
$ws.Range("B16").Value      = "GP01"   # was GP1
$ws.Range("B17:B18").Value  = "GP02"   # was GP2
$ws.Range("B19:B20").Value  = "GP03"   # was GP3
$ws.Range("B21:B25").Value  = "BP01"   # was BP1
$ws.Range("B26:B30").Value  = "BP02"   # was BP2
$ws.Range("B31:B60").Value  = "BP03"   # was BP3

# --- reset view / selection back to the top of the sheet ---
[void]$ws.Range("A1").Select()
